# Update the "Förändrad" date column (C) for every data row (2-135)
# from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C135").Value = 45184
